# Fix wording on slide 2 ("Idee"): the invalid-input bullet point should say
# that the input is not *accepted* ("akzeptiert") instead of not *saved*
# ("gespeichert").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$old = "ist der Input ungültig: wird die Schrift rot und die Eingabe wird nicht gespeichert "
$new = "ist der Input ungültig: wird die Schrift rot und die Eingabe wird nicht akzeptiert "

$full = $tr.Text
$idx  = $full.IndexOf($old)

if ($idx -ge 0) {
    # Characters() is 1-based, IndexOf() is 0-based.
    $sub = $tr.Characters($idx + 1, $old.Length)
    $sub.Text = $new
} else {
    # Fallback: text already matches the target, or structure differs -
    # try replacing the whole second paragraph directly.
    $para = $tr.Paragraphs(2)
    if ($para.Text -eq "ist der Input ungültig: wird die Schrift rot und die Eingabe wird nicht gespeichert") {
        $para.Text = "ist der Input ungültig: wird die Schrift rot und die Eingabe wird nicht akzeptiert "
    }
}
